$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.693996120467564
$ws.Range("K2").Value = 0.691681993781983
$ws.Range("L2").Value = 0.754855071318372
$ws.Range("N2").Value = 0.622335501686223

$ws.Range("B3").Value = 0.641188431291775
$ws.Range("K3").Value = 0.550394111458431
$ws.Range("L3").Value = 0.756622372880961
$ws.Range("N3").Value = 0.57421952327369

$ws.Range("B4").Value = 0.704965329416964
$ws.Range("K4").Value = 0.703727559204962
$ws.Range("L4").Value = 0.778963825426238
$ws.Range("N4").Value = 0.561007628053552

$ws.Range("B5").Value = 0.451810364536854
$ws.Range("N5").Value = 0.401158464849325

$ws.Range("B6").Value = 0.592122368373113
$ws.Range("K6").Value = 0.584137117430032
$ws.Range("L6").Value = 0.573182308826315
$ws.Range("N6").Value = 0.554219568199052
